$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.195.10"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "'3.148.43"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D5").Value = "'535.20"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "'138.78"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'3.147.58"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.468"
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "'0.414"
$ws.Range("E12").Value = "  +4.28%  "
$ws.Range("D13").Value = "'3.686.86"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'25.72"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'58.274.66"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'3.164.05"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'12.71"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("E21").Value = "  +2.90%  "
$ws.Range("D22").Value = "'360.47"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'69.18"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'0.506"
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'0.0₃0882"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").Value = "'159.03"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'6.09"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "'25.99"
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "'1.69"
$ws.Range("E39").Value = "  +4.48%  "
$ws.Range("D40").Value = "'0.0670"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'2.513.17"
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("D42").Value = "'0.704"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "'37.46"
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "'3.190.19"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "'19.84"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("E51").Value = "  -4.03%  "
